# Feature/wmt 414/display arms assessments
# Add an "ARMS" tab to the workbook with the Assessment header row,
# formatted to match the other header rows in the workbook.

$wb = $excel.ActiveWorkbook

# --- Add the new "ARMS" worksheet as the last sheet (after "CMS") ---
$cmsSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cmsSheet)
$ws.Name = "ARMS"

# --- Header values (row 1) ---
$ws.Range("A1").Value = "Assessment_Staff_Name"
$ws.Range("B1").Value = "Assessment_Staff_Key"
$ws.Range("C1").Value = "Assessment_Staff_Grade"
$ws.Range("D1").Value = "Assessmentent_Team_Key"
$ws.Range("E1").Value = "Assessment_Provider_Code"
$ws.Range("F1").Value = "CRN"
$ws.Range("G1").Value = "Disposal_or_Release_Date"
$ws.Range("I1").Value = "SO_Registration_Date"
$ws.Range("H1").Value = "Sentence_Type"

# --- Row height to match the other "header row" sheets ---
$ws.Range("A1:I1").RowHeight = 37

# --- Common header formatting: bold white 9pt Arial on a blue fill ---
$headerRange = $ws.Range("A1:I1")
$headerRange.WrapText = $true
$headerRange.HorizontalAlignment = -4131
$headerFont = $headerRange.Font
$headerFont.Name = "Arial"
$headerFont.Bold = $true
$headerFont.Size = 9
$headerFont.Color = 16777215
$headerRange.Interior.Color = 10511371

# --- Borders ---
# A1: full box (left/right/top thin blue, bottom thin grey)
$a1 = $ws.Range("A1")
$a1.Borders.Item(7).LineStyle = 1
$a1.Borders.Item(7).Color = 10909496
$a1.Borders.Item(10).LineStyle = 1
$a1.Borders.Item(10).Color = 10909496
$a1.Borders.Item(8).LineStyle = 1
$a1.Borders.Item(8).Color = 10909496
$a1.Borders.Item(9).LineStyle = 1
$a1.Borders.Item(9).Color = 11642277

# B1:G1, I1: right + top thin blue, bottom thin grey (no left border)
$midRange = $ws.Range("B1:G1")
$midRange.Borders.Item(10).LineStyle = 1
$midRange.Borders.Item(10).Color = 10909496
$midRange.Borders.Item(8).LineStyle = 1
$midRange.Borders.Item(8).Color = 10909496
$midRange.Borders.Item(9).LineStyle = 1
$midRange.Borders.Item(9).Color = 11642277

# C1 additionally uses a text number format and no wrap
$c1 = $ws.Range("C1")
$c1.NumberFormat = "@"
$c1.WrapText = $false

# H1 and I1: only a thin blue right border (no top/bottom)
$endRange = $ws.Range("H1:I1")
$endRange.Borders.Item(10).LineStyle = 1
$endRange.Borders.Item(10).Color = 10909496

$i1 = $ws.Range("I1")
$i1.Borders.Item(10).LineStyle = 1
$i1.Borders.Item(10).Color = 10909496

# --- Selection / active cell on the new sheet ---
$ws.Range("H1").Select()

Write-Output "ARMS sheet added"
